$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 22:52"

# Tenerife row (row 33): Casos activos, Recuperados, Muertes
$ws.Range("C33").Value = 337
$ws.Range("D33").Value = 829
$ws.Range("E33").Value = 71

# Gran Canaria row (row 50): Casos activos, Recuperados
$ws.Range("C50").Value = 193
$ws.Range("D50").Value = 236
